# The sheet tracks weekly price observations for "Sandia" (watermelon) at
# "Vega Modelo de Temuco". This commit adds a new weekly record, inserted
# as row 287 (pushing every existing row from 287 on down by one, and
# extending the used range from R324 to R325).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 287; everything currently at/after 287 shifts down
# by one (old 287 -> 288, ..., old 324 -> 325), which matches the target
# workbook exactly.
$ws.Rows.Item(287).Insert()

# Populate the newly inserted row 287 with the new weekly observation.
$ws.Range("A287").Value = 10
$ws.Range("B287").Value = "Vega Modelo de Temuco"
$ws.Range("C287").Value = "La Araucanía"
$ws.Range("D287").Value = 44491
$ws.Range("E287").Value = 9
$ws.Range("F287").Value = 100112028
$ws.Range("G287").Value = "Sandia"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Segunda"
$ws.Range("J287").Value = 300
$ws.Range("K287").Value = 1000
$ws.Range("L287").Value = 1000
$ws.Range("M287").Value = 1000
$ws.Range("N287").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O287").Value = "Perú"
$ws.Range("P287").Value = 1000
$ws.Range("Q287").Value = 1
$ws.Range("R287").Value = "Hortaliza"
